$d = $word.ActiveDocument

# Locate the paragraph that reads "Core: TBX-Core" and strip the
# stray "TBX-" so it reads "Core: Core" (module name fix).
$found = $false
foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    if ($pr.Text -like "Core: TBX-Core*") {
        $paraStart = $pr.Start
        $text = $pr.Text
        $offset = $text.IndexOf("TBX-")
        $delStart = $paraStart + $offset
        $delEnd = $delStart + 4

        # Remove the "TBX-" substring in place.
        $d.Range($delStart, $delEnd).Delete()

        # Word re-anchors its hidden "_GoBack" bookmark (last-edit marker)
        # at the location of the most recent edit; adding a bookmark with
        # that reserved name automatically replaces any prior occurrence.
        $editRange = $d.Range($delStart, $delStart)
        $d.Bookmarks.Add("_GoBack", $editRange)

        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not find the 'Core: TBX-Core' paragraph to fix"
}

$d.Save()
